$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "68.520.33"
Set-TextValue "E2" "  +1.72%  "
Set-TextValue "D3" "3.761.20"
Set-TextValue "E3" "  +0.26%  "
Set-TextValue "E4" "  +0.13%  "
Set-TextValue "D5" "596.54"
Set-TextValue "E5" "  +0.40%  "
Set-TextValue "D6" "168.40"
Set-TextValue "E6" "  -0.74%  "
Set-TextValue "D7" "3.758.92"
Set-TextValue "E7" "  +0.29%  "
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "D9" "0.522"
Set-TextValue "E9" "  -0.42%  "
Set-TextValue "D10" "0.161"
Set-TextValue "E10" "  -2.54%  "
Set-TextValue "D11" "6.49"
Set-TextValue "E11" "  +0.36%  "
Set-TextValue "D12" "0.451"
Set-TextValue "E12" "  -0.44%  "
Set-TextValue "D13" "0.0000261"
Set-TextValue "E13" "  -5.07%  "
Set-TextValue "D14" "36.55"
Set-TextValue "E14" "  +0.29%  "
Set-TextValue "D15" "4.395.74"
Set-TextValue "E15" "  +0.24%  "
Set-TextValue "D16" "3.760.44"
Set-TextValue "E16" "  -0.06%  "
Set-TextValue "D17" "68.564.93"
Set-TextValue "E17" "  +1.76%  "
Set-TextValue "D18" "18.06"
Set-TextValue "E18" "  -2.60%  "
Set-TextValue "E19" "  -1.87%  "
Set-TextValue "E20" "  -0.41%  "
Set-TextValue "D21" "10.88"
Set-TextValue "E21" "  +3.53%  "
Set-TextValue "D22" "468.33"
Set-TextValue "E22" "  +0.30%  "
Set-TextValue "D23" "0.704"
Set-TextValue "E23" "  -1.86%  "
Set-TextValue "E24" "  +0.86%  "
Set-TextValue "E25" "  -0.25%  "
Set-TextValue "E26" "  +0.49%  "
Set-TextValue "D27" "12.13"
Set-TextValue "E27" "  +0.05%  "
Set-TextValue "D28" "10.19"
Set-TextValue "E28" "  -1.00%  "
Set-TextValue "E29" "  +0.04%  "
Set-TextValue "D30" "3.911.84"
Set-TextValue "E30" "  +0.13%  "
Set-TextValue "D31" "2.80"
Set-TextValue "E31" "  -3.16%  "
Set-TextValue "D32" "7.39"
Set-TextValue "E32" "  -3.15%  "
Set-TextValue "D33" "30.10"
Set-TextValue "E33" "  -1.32%  "
Set-TextValue "D34" "2.21"
Set-TextValue "E34" "  -0.60%  "
Set-TextValue "D35" "9.33"
Set-TextValue "E35" "  +2.47%  "
Set-TextValue "D36" "1.00"
Set-TextValue "D37" "3.718.48"
Set-TextValue "E37" "  +0.03%  "
Set-TextValue "E38" "  -1.75%  "
Set-TextValue "D39" "3.46"
Set-TextValue "E39" "  -9.11%  "
Set-TextValue "E40" "  +1.39%  "
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  +0.67%  "
Set-TextValue "D42" "5.84"
Set-TextValue "E42" "  +0.25%  "
Set-TextValue "E43" "  +0.26%  "
Set-TextValue "E45" "  -0.99%  "
Set-TextValue "D46" "1.96"
Set-TextValue "E46" "  +1.20%  "
Set-TextValue "D47" "43.64"
Set-TextValue "E47" "  +11.94%  "
Set-TextValue "D48" "8.62"
Set-TextValue "E48" "  -0.99%  "
Set-TextValue "D49" "45.97"
Set-TextValue "E49" "  +0.35%  "
Set-TextValue "D50" "397.91"
Set-TextValue "E50" "  +0.37%  "
Set-TextValue "D51" "146.29"
Set-TextValue "E51" "  +6.16%  "
